$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.401.65"
$ws.Range("E2").Value = "  -1.69%  "

# Row 3
$ws.Range("D3").Value = "2.619.25"
$ws.Range("E3").Value = "  +0.76%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.98"
$ws.Range("E5").Value = "  -0.86%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.11"
$ws.Range("E6").Value = "  +0.69%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").Value = "  +0.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.92"
$ws.Range("E9").Value = "  +6.68%  "

# Row 10
$ws.Range("E10").Value = "  -1.96%  "

# Row 11
$ws.Range("E11").Value = "  -0.34%  "

# Row 12
$ws.Range("E12").Value = "  +0.96%  "

# Row 13
$ws.Range("D13").Value = "3.087.85"
$ws.Range("E13").Value = "  +1.00%  "

# Row 14
$ws.Range("D14").Value = "58.338.41"
$ws.Range("E14").Value = "  -1.67%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.64"
$ws.Range("E15").Value = "  -0.18%  "

# Row 16
$ws.Range("D16").Value = "2.622.87"
$ws.Range("E16").Value = "  +0.98%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  -1.32%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.38"
$ws.Range("E18").Value = "  +0.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "333.70"
$ws.Range("E19").Value = "  -2.01%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.09"
$ws.Range("E20").Value = "  +0.18%  "

# Row 21
$ws.Range("E21").Value = "  -2.18%  "

# Row 22
$ws.Range("E22").Value = "  +0.01%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.33"
$ws.Range("E23").Value = "  -1.53%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.414"
$ws.Range("E24").Value = "  +1.50%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.02%  "

# Row 26
$ws.Range("E26").Value = "  -1.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.06"
$ws.Range("E27").Value = "  -2.01%  "

# Row 28
$ws.Range("B28").Value = "USDe"
$ws.Range("C28").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.04%  "

# Row 29
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0731"
$ws.Range("E29").Value = "  -1.46%  "

# Row 30
$ws.Range("E30").Value = "  -1.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.84"
$ws.Range("E31").Value = "  +0.47%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.73"
$ws.Range("E32").Value = "  -0.21%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.01"
$ws.Range("E33").Value = "  +0.15%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.87"
$ws.Range("E34").Value = "  -2.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.847"
$ws.Range("E35").Value = "  +0.82%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.09"
$ws.Range("E36").Value = "  -1.66%  "

# Row 37
$ws.Range("E37").Value = "  -3.68%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.805"
$ws.Range("E38").Value = "  -1.94%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.55"
$ws.Range("E39").Value = "  +0.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "279.12"

# Row 41
$ws.Range("E41").Value = "  +0.05%  "

# Row 42
$ws.Range("E42").Value = "  -0.86%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.67"
$ws.Range("E43").Value = "  -0.51%  "

# Row 44
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0527"
$ws.Range("E44").Value = "  +0.61%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.90"
$ws.Range("E45").Value = "  +2.45%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0934"
$ws.Range("E46").Value = "  -1.95%  "

# Row 47
$ws.Range("E47").Value = "  +0.23%  "

# Row 48
$ws.Range("D48").Value = "1.937.08"
$ws.Range("E48").Value = "  -0.12%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.41"
$ws.Range("E49").Value = "  -1.68%  "

# Row 50
$ws.Range("E50").Value = "  -4.02%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.55"
$ws.Range("E51").Value = "  +1.20%  "
